$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

$aboutSheet.Range("A2").Value = "Version: $newVersion"
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Shakhtinskaya Coal Mine, Kazakhstan, M1440, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

for ($row = 2; $row -le 9; $row++) {
    $dataSheet.Range("S$row").Value = $newVersion
}
